$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "BTC"
$ws.Cells.Item(2, 3).Value = "Bitcoin"
$ws.Cells.Item(2, 4).Value = 25802
$ws.Cells.Item(2, 5).Value = 500595488301
$ws.Cells.Item(2, 6).Value = 6369370944
$ws.Cells.Item(2, 7).Value = 0.44327

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "ETH"
$ws.Cells.Item(3, 3).Value = "Ethereum"
$ws.Cells.Item(3, 4).Value = 1751.95
$ws.Cells.Item(3, 5).Value = 210696959252
$ws.Cells.Item(3, 6).Value = 6283205999
$ws.Cells.Item(3, 7).Value = 0.42571

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "USDT"
$ws.Cells.Item(4, 3).Value = "Tether"
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 83408115139
$ws.Cells.Item(4, 6).Value = 20080085200
$ws.Cells.Item(4, 7).Value = -0.00122

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "BNB"
$ws.Cells.Item(5, 3).Value = "BNB"
$ws.Cells.Item(5, 4).Value = 236.62
$ws.Cells.Item(5, 5).Value = 36866596856
$ws.Cells.Item(5, 6).Value = 743668096
$ws.Cells.Item(5, 7).Value = -0.60977

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "USDC"
$ws.Cells.Item(6, 3).Value = "USD Coin"
$ws.Cells.Item(6, 4).Value = 0.999848
$ws.Cells.Item(6, 5).Value = 28371716689
$ws.Cells.Item(6, 6).Value = 2471803424
$ws.Cells.Item(6, 7).Value = -0.02477

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "XRP"
$ws.Cells.Item(7, 3).Value = "XRP"
$ws.Cells.Item(7, 4).Value = 0.511595
$ws.Cells.Item(7, 5).Value = 26556696804
$ws.Cells.Item(7, 6).Value = 1060243833
$ws.Cells.Item(7, 7).Value = 3.22591

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "STETH"
$ws.Cells.Item(8, 3).Value = "Lido Staked Ether"
$ws.Cells.Item(8, 4).Value = 1752.19
$ws.Cells.Item(8, 5).Value = 12562695228
$ws.Cells.Item(8, 6).Value = 23226503
$ws.Cells.Item(8, 7).Value = 0.34485

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "ADA"
$ws.Cells.Item(9, 3).Value = "Cardano"
$ws.Cells.Item(9, 4).Value = 0.270096
$ws.Cells.Item(9, 5).Value = 9444479618
$ws.Cells.Item(9, 6).Value = 733563422
$ws.Cells.Item(9, 7).Value = 8.14471

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "DOGE"
$ws.Cells.Item(10, 3).Value = "Dogecoin"
$ws.Cells.Item(10, 4).Value = 0.062276
$ws.Cells.Item(10, 5).Value = 8697784214
$ws.Cells.Item(10, 6).Value = 491575595
$ws.Cells.Item(10, 7).Value = 4.22588

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "TRX"
$ws.Cells.Item(11, 3).Value = "TRON"
$ws.Cells.Item(11, 4).Value = 0.069215
$ws.Cells.Item(11, 5).Value = 6236833750
$ws.Cells.Item(11, 6).Value = 297841168
$ws.Cells.Item(11, 7).Value = 2.07891

$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "SOL"
$ws.Cells.Item(12, 3).Value = "Solana"
$ws.Cells.Item(12, 4).Value = 15.57
$ws.Cells.Item(12, 5).Value = 6195433378
$ws.Cells.Item(12, 6).Value = 709387933
$ws.Cells.Item(12, 7).Value = 5.4866

$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "MATIC"
$ws.Cells.Item(13, 3).Value = "Polygon"
$ws.Cells.Item(13, 4).Value = 0.627524
$ws.Cells.Item(13, 5).Value = 5824127171
$ws.Cells.Item(13, 6).Value = 679709510
$ws.Cells.Item(13, 7).Value = 7.89528

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "LTC"
$ws.Cells.Item(14, 3).Value = "Litecoin"
$ws.Cells.Item(14, 4).Value = 78.21
$ws.Cells.Item(14, 5).Value = 5721375955
$ws.Cells.Item(14, 6).Value = 657407635
$ws.Cells.Item(14, 7).Value = 1.15937

$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "DOT"
$ws.Cells.Item(15, 3).Value = "Polkadot"
$ws.Cells.Item(15, 4).Value = 4.49
$ws.Cells.Item(15, 5).Value = 5580258420
$ws.Cells.Item(15, 6).Value = 151747287
$ws.Cells.Item(15, 7).Value = 0.49926

$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "BUSD"
$ws.Cells.Item(16, 3).Value = "Binance USD"
$ws.Cells.Item(16, 4).Value = 0.999738
$ws.Cells.Item(16, 5).Value = 4779419798
$ws.Cells.Item(16, 6).Value = 1122823868
$ws.Cells.Item(16, 7).Value = -0.07326

$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "DAI"
$ws.Cells.Item(17, 3).Value = "Dai"
$ws.Cells.Item(17, 4).Value = 0.999956
$ws.Cells.Item(17, 5).Value = 4542725456
$ws.Cells.Item(17, 6).Value = 102850792
$ws.Cells.Item(17, 7).Value = 0.03421

$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "WBTC"
$ws.Cells.Item(18, 3).Value = "Wrapped Bitcoin"
$ws.Cells.Item(18, 4).Value = 25791
$ws.Cells.Item(18, 5).Value = 4042830718
$ws.Cells.Item(18, 6).Value = 86871186
$ws.Cells.Item(18, 7).Value = 0.25331

$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "AVAX"
$ws.Cells.Item(19, 3).Value = "Avalanche"
$ws.Cells.Item(19, 4).Value = 11.7
$ws.Cells.Item(19, 5).Value = 4031097052
$ws.Cells.Item(19, 6).Value = 217402598
$ws.Cells.Item(19, 7).Value = 0.87344

$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "SHIB"
$ws.Cells.Item(20, 3).Value = "Shiba Inu"
$ws.Cells.Item(20, 4).Value = 0.00000673
$ws.Cells.Item(20, 5).Value = 3966292051
$ws.Cells.Item(20, 6).Value = 162687111
$ws.Cells.Item(20, 7).Value = 2.52758

$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "LEO"
$ws.Cells.Item(21, 3).Value = "LEO Token"
$ws.Cells.Item(21, 4).Value = 3.55
$ws.Cells.Item(21, 5).Value = 3293095380
$ws.Cells.Item(21, 6).Value = 638139
$ws.Cells.Item(21, 7).Value = 0.89665

$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "UNI"
$ws.Cells.Item(22, 3).Value = "Uniswap"
$ws.Cells.Item(22, 4).Value = 4.08
$ws.Cells.Item(22, 5).Value = 3073134648
$ws.Cells.Item(22, 6).Value = 57263270
$ws.Cells.Item(22, 7).Value = 2.11153

$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "LINK"
$ws.Cells.Item(23, 3).Value = "Chainlink"
$ws.Cells.Item(23, 4).Value = 5.19
$ws.Cells.Item(23, 5).Value = 2680024536
$ws.Cells.Item(23, 6).Value = 247445264
$ws.Cells.Item(23, 7).Value = 2.26608

$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "XMR"
$ws.Cells.Item(24, 3).Value = "Monero"
$ws.Cells.Item(24, 4).Value = 136.87
$ws.Cells.Item(24, 5).Value = 2486158995
$ws.Cells.Item(24, 6).Value = 54949285
$ws.Cells.Item(24, 7).Value = -0.71641

$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = "OKB"
$ws.Cells.Item(25, 3).Value = "OKB"
$ws.Cells.Item(25, 4).Value = 40.85
$ws.Cells.Item(25, 5).Value = 2455124330
$ws.Cells.Item(25, 6).Value = 9052172
$ws.Cells.Item(25, 7).Value = -1.70293

$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "ATOM"
$ws.Cells.Item(26, 3).Value = "Cosmos Hub"
$ws.Cells.Item(26, 4).Value = 8.26
$ws.Cells.Item(26, 5).Value = 2416569452
$ws.Cells.Item(26, 6).Value = 99742775
$ws.Cells.Item(26, 7).Value = 4.18942

$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "XLM"
$ws.Cells.Item(27, 3).Value = "Stellar"
$ws.Cells.Item(27, 4).Value = 0.082884
$ws.Cells.Item(27, 5).Value = 2226753924
$ws.Cells.Item(27, 6).Value = 54139534
$ws.Cells.Item(27, 7).Value = 2.23928

$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "ETC"
$ws.Cells.Item(28, 3).Value = "Ethereum Classic"
$ws.Cells.Item(28, 4).Value = 15.22
$ws.Cells.Item(28, 5).Value = 2155780934
$ws.Cells.Item(28, 6).Value = 104106505
$ws.Cells.Item(28, 7).Value = 4.2874

$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = "TON"
$ws.Cells.Item(29, 3).Value = "Toncoin"
$ws.Cells.Item(29, 4).Value = 1.46
$ws.Cells.Item(29, 5).Value = 2153693752
$ws.Cells.Item(29, 6).Value = 7730983
$ws.Cells.Item(29, 7).Value = -2.13301

$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = "TUSD"
$ws.Cells.Item(30, 3).Value = "TrueUSD"
$ws.Cells.Item(30, 4).Value = 0.998542
$ws.Cells.Item(30, 5).Value = 2038588708
$ws.Cells.Item(30, 6).Value = 1621114950
$ws.Cells.Item(30, 7).Value = 0.21368

$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = "BCH"
$ws.Cells.Item(31, 3).Value = "Bitcoin Cash"
$ws.Cells.Item(31, 4).Value = 102.71
$ws.Cells.Item(31, 5).Value = 1993609792
$ws.Cells.Item(31, 6).Value = 73106056
$ws.Cells.Item(31, 7).Value = 1.54678

$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = "ICP"
$ws.Cells.Item(32, 3).Value = "Internet Computer"
$ws.Cells.Item(32, 4).Value = 3.75
$ws.Cells.Item(32, 5).Value = 1631291793
$ws.Cells.Item(32, 6).Value = 31307051
$ws.Cells.Item(32, 7).Value = -0.5126

$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = "LDO"
$ws.Cells.Item(33, 3).Value = "Lido DAO"
$ws.Cells.Item(33, 4).Value = 1.78
$ws.Cells.Item(33, 5).Value = 1565492493
$ws.Cells.Item(33, 6).Value = 49254589
$ws.Cells.Item(33, 7).Value = -3.63985

$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = "QNT"
$ws.Cells.Item(34, 3).Value = "Quant"
$ws.Cells.Item(34, 4).Value = 101.98
$ws.Cells.Item(34, 5).Value = 1483065603
$ws.Cells.Item(34, 6).Value = 19896426
$ws.Cells.Item(34, 7).Value = -2.06121

$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = "FIL"
$ws.Cells.Item(35, 3).Value = "Filecoin"
$ws.Cells.Item(35, 4).Value = 3.43
$ws.Cells.Item(35, 5).Value = 1472491766
$ws.Cells.Item(35, 6).Value = 173858361
$ws.Cells.Item(35, 7).Value = 2.39789

$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = "CRO"
$ws.Cells.Item(36, 3).Value = "Cronos"
$ws.Cells.Item(36, 4).Value = 0.055226
$ws.Cells.Item(36, 5).Value = 1441425857
$ws.Cells.Item(36, 6).Value = 10155746
$ws.Cells.Item(36, 7).Value = 7.32196

$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "HBAR"
$ws.Cells.Item(37, 3).Value = "Hedera"
$ws.Cells.Item(37, 4).Value = 0.04421969
$ws.Cells.Item(37, 5).Value = 1397971881
$ws.Cells.Item(37, 6).Value = 32576736
$ws.Cells.Item(37, 7).Value = -0.49273

$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "ARB"
$ws.Cells.Item(38, 3).Value = "Arbitrum"
$ws.Cells.Item(38, 4).Value = 1.003
$ws.Cells.Item(38, 5).Value = 1279626021
$ws.Cells.Item(38, 6).Value = 306994497
$ws.Cells.Item(38, 7).Value = -1.27178

$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = "APT"
$ws.Cells.Item(39, 3).Value = "Aptos"
$ws.Cells.Item(39, 4).Value = 6
$ws.Cells.Item(39, 5).Value = 1203260361
$ws.Cells.Item(39, 6).Value = 88383069
$ws.Cells.Item(39, 7).Value = 0.51843

$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = "VET"
$ws.Cells.Item(40, 3).Value = "VeChain"
$ws.Cells.Item(40, 4).Value = 0.01571441
$ws.Cells.Item(40, 5).Value = 1142103912
$ws.Cells.Item(40, 6).Value = 51790710
$ws.Cells.Item(40, 7).Value = 4.69743

$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = "NEAR"
$ws.Cells.Item(41, 3).Value = "NEAR Protocol"
$ws.Cells.Item(41, 4).Value = 1.21
$ws.Cells.Item(41, 5).Value = 1104636075
$ws.Cells.Item(41, 6).Value = 72960357
$ws.Cells.Item(41, 7).Value = 0.47387

$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = "USDP"
$ws.Cells.Item(42, 3).Value = "Pax Dollar"
$ws.Cells.Item(42, 4).Value = 0.998169
$ws.Cells.Item(42, 5).Value = 1004562920
$ws.Cells.Item(42, 6).Value = 1659208
$ws.Cells.Item(42, 7).Value = -0.20889

$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "FRAX"
$ws.Cells.Item(43, 3).Value = "Frax"
$ws.Cells.Item(43, 4).Value = 0.998468
$ws.Cells.Item(43, 5).Value = 1002513369
$ws.Cells.Item(43, 6).Value = 8013847
$ws.Cells.Item(43, 7).Value = -0.12378

$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "GRT"
$ws.Cells.Item(44, 3).Value = "The Graph"
$ws.Cells.Item(44, 4).Value = 0.100782
$ws.Cells.Item(44, 5).Value = 907155639
$ws.Cells.Item(44, 6).Value = 56806192
$ws.Cells.Item(44, 7).Value = 4.2034

$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = "BSCX"
$ws.Cells.Item(45, 3).Value = "BSCEX"
$ws.Cells.Item(45, 4).Value = 237.12
$ws.Cells.Item(45, 5).Value = 904393606
$ws.Cells.Item(45, 6).Value = 1239178
$ws.Cells.Item(45, 7).Value = -0.5628

$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "RPL"
$ws.Cells.Item(46, 3).Value = "Rocket Pool"
$ws.Cells.Item(46, 4).Value = 44.22
$ws.Cells.Item(46, 5).Value = 861688515
$ws.Cells.Item(46, 6).Value = 3460330
$ws.Cells.Item(46, 7).Value = -0.13388

$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "APE"
$ws.Cells.Item(47, 3).Value = "ApeCoin"
$ws.Cells.Item(47, 4).Value = 2.33
$ws.Cells.Item(47, 5).Value = 855871892
$ws.Cells.Item(47, 6).Value = 170097834
$ws.Cells.Item(47, 7).Value = 0.02569

$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "RETH"
$ws.Cells.Item(48, 3).Value = "Rocket Pool ETH"
$ws.Cells.Item(48, 4).Value = 1882.3
$ws.Cells.Item(48, 5).Value = 829056230
$ws.Cells.Item(48, 6).Value = 2109249
$ws.Cells.Item(48, 7).Value = 0.30041

$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = "ALGO"
$ws.Cells.Item(49, 3).Value = "Algorand"
$ws.Cells.Item(49, 4).Value = 0.110078
$ws.Cells.Item(49, 5).Value = 796730452
$ws.Cells.Item(49, 6).Value = 51655433
$ws.Cells.Item(49, 7).Value = 1.36527

$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = "EGLD"
$ws.Cells.Item(50, 3).Value = "MultiversX"
$ws.Cells.Item(50, 4).Value = 30.2
$ws.Cells.Item(50, 5).Value = 770327338
$ws.Cells.Item(50, 6).Value = 14143998
$ws.Cells.Item(50, 7).Value = -0.1087

$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = "STX"
$ws.Cells.Item(51, 3).Value = "Stacks"
$ws.Cells.Item(51, 4).Value = 0.552745
$ws.Cells.Item(51, 5).Value = 762591631
$ws.Cells.Item(51, 6).Value = 55596251
$ws.Cells.Item(51, 7).Value = 3.76432
